$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") — copy formatting from the
# existing header cell H1 so they pick up the same bold/border/centered
# style (style index 1), then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for the new I (I0) and J (IF) columns, rows 2-22.
$data = @(
    @(1, 6),
    @(1, 6),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 7),
    @(1, 6),
    @(1, 6),
    @(1, 7),
    @(1, 5),
    @(1, 7),
    @(1, 5),
    @(7, 8),
    @(1, 3),
    @(5, 8),
    @(6, 7),
    @(1, 3),
    @(5, 7),
    @(1, 4),
    @(4, 6),
    @(3, 4)
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
    $r++
}
